# Update Main.xlsx "Rules" sheet as part of the "update file with jgit" commit:
#  - cell E8 changes from "Good Morning" to "GIT UPDATE"
#  - the active selection ends up on E8

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E8").Value = "GIT UPDATE"
$ws.Range("E8").Select()
